$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-06 Friday" "2024-09-07 Saturday"

Replace-Text "90÷2=" "79÷7="
Replace-Text "83÷8=" "33÷4="
Replace-Text "91÷2=" "74÷9="
Replace-Text "85÷6=" "34÷8="
Replace-Text "35÷3=" "32÷9="

Replace-Text "62÷5=" "34÷7="
Replace-Text "52÷4=" "58÷2="
Replace-Text "19÷4=" "60÷8="
Replace-Text "13÷7=" "63÷5="
Replace-Text "37÷3=" "22÷7="

Replace-Text "50÷7=" "73÷7="
Replace-Text "16÷9=" "98÷7="
Replace-Text "89÷2=" "79÷3="
Replace-Text "97÷9=" "65÷2="
Replace-Text "43÷5=" "46÷4="

Replace-Text "12÷2=" "39÷2="
Replace-Text "61÷7=" "95÷3="
Replace-Text "51÷7=" "11÷9="
Replace-Text "64÷3=" "35÷8="
Replace-Text "86÷2=" "42÷6="

Replace-Text "49÷6=" "11÷4="
Replace-Text "97÷8=" "30÷9="
Replace-Text "47÷2=" "26÷8="
Replace-Text "67÷7=" "80÷2="
Replace-Text "98÷2=" "23÷6="
